# Insert a new data row at row 23 (pushes existing rows 23:95 down to 24:96)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(23).Insert()

# Populate the newly inserted row 23 with its data
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = 44672
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100102
$ws.Range("H23").Value = "Cítricos"
$ws.Range("I23").Value = 100102004
$ws.Range("J23").Value = "Mandarina"
$ws.Range("K23").Value = "Murcott"
$ws.Range("L23").Value = "Tercera"
$ws.Range("M23").Value = 300
$ws.Range("N23").Value = 13000
$ws.Range("O23").Value = 14000
$ws.Range("P23").Value = 13500
$ws.Range("Q23").Value = "$/caja 20 kilos"
$ws.Range("R23").Value = "Región de Coquimbo"
$ws.Range("S23").Value = 675
$ws.Range("T23").Value = 20
